$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NSE:ADFFOODS"
$ws.Range("C2").Value = "NSE:AARVI"
$ws.Range("F2").Value = "NSE:BAJAJFINSV"

$ws.Range("B3").Value = "NSE:ALICON"
$ws.Range("C3").Value = "NSE:ACI"
$ws.Range("D3").Value = "NSE:IPCALAB"
$ws.Range("F3").Value = "NSE:ICICIBANK"

$ws.Range("B4").Value = "NSE:APOLSINHOT"
$ws.Range("C4").Value = "NSE:BCLIND"
$ws.Range("D4").Value = "NSE:JKCEMENT"
$ws.Range("F4").ClearContents()

$ws.Range("B5").Value = "NSE:BAJAJFINSV"
$ws.Range("C5").Value = "NSE:CDSL"

$ws.Range("B6").Value = "NSE:CENTURYTEX"
$ws.Range("C6").Value = "NSE:DIAMONDYD"

$ws.Range("B7").Value = "NSE:DHANUKA"
$ws.Range("C7").Value = "NSE:FAZE3Q"

$ws.Range("B8").Value = "NSE:GLOBAL"
$ws.Range("C8").Value = "NSE:FMGOETZE"

$ws.Range("B9").Value = "NSE:GOKEX"
$ws.Range("C9").Value = "NSE:FOODSIN"

$ws.Range("B10").Value = "NSE:IRISDOREME"
$ws.Range("C10").Value = "NSE:HPL"

$ws.Range("B11").Value = "NSE:NINSYS"
$ws.Range("C11").Value = "NSE:IFCI"

$ws.Range("B12").Value = "NSE:PAVNAIND"
$ws.Range("C12").Value = "NSE:KEI"

$ws.Range("B13").Value = "NSE:RAIN"
$ws.Range("C13").Value = "NSE:KRITI"

$ws.Range("B14").Value = "NSE:SAMBHAAV"
$ws.Range("C14").Value = "NSE:KSCL"

$ws.Range("B15").ClearContents()
$ws.Range("C15").Value = "NSE:LAMBODHARA"

$ws.Range("B16").ClearContents()
$ws.Range("C16").Value = "NSE:NAVNETEDUL"

$ws.Range("B17").ClearContents()
$ws.Range("C17").Value = "NSE:NESCO"

$ws.Range("B18").ClearContents()
$ws.Range("C18").Value = "NSE:ORIENTBELL"

$ws.Range("B19").ClearContents()
$ws.Range("C19").Value = "NSE:PRITI"

$ws.Range("B20").ClearContents()
$ws.Range("C20").Value = "NSE:RGL"

$ws.Range("B21").ClearContents()
$ws.Range("C21").Value = "NSE:RHFL"

$ws.Range("B22").ClearContents()
$ws.Range("C22").Value = "NSE:RITCO"

$ws.Range("B23").ClearContents()
$ws.Range("C23").Value = "NSE:RML"

$ws.Range("B24").ClearContents()
$ws.Range("C24").Value = "NSE:RPOWER"

$ws.Range("B25").ClearContents()
$ws.Range("C25").Value = "NSE:SAKHTISUG"

$ws.Range("B26").ClearContents()
$ws.Range("C26").Value = "NSE:SALASAR"

# Remove rows 27-29 (bottom-up so row indices stay valid)
$ws.Rows.Item(29).Delete()
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(27).Delete()
